$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 187: update existing row with new/changed values
$ws.Range("B187").Value = 33970.5
$ws.Range("C187").Value = 14713.9
$ws.Range("D187").Value = 6903.9
$ws.Range("F187").Value = 6455.8
$ws.Range("G187").Value = 15132.1
$ws.Range("J187").Value = 1527.9
$ws.Range("K187").Value = 69144.4
$ws.Range("L187").Value = 3969.6
$ws.Range("M187").Value = 1603.1
$ws.Range("O187").Value = 1391.9
$ws.Range("P187").Value = 73635.2
$ws.Range("Q187").Value = 108843.7
$ws.Range("R187").Value = 1306.6
$ws.Range("S187").Value = 50557.7
$ws.Range("T187").Value = 17747

# Row 188: new row for 21-09-2021
$ws.Range("A188").Value = "21-09-2021"
$ws.Range("B188").Value = 33919.8
$ws.Range("C188").Value = 14746.4
$ws.Range("D188").Value = 6981
$ws.Range("E188").Value = 29839.7
$ws.Range("F188").Value = 6552.7
$ws.Range("G188").Value = 15348.5
$ws.Range("J188").Value = 1530.4
$ws.Range("K188").Value = 69420
$ws.Range("L188").Value = 3970.5
$ws.Range("M188").Value = 1614.9
$ws.Range("O188").Value = 1385.6
$ws.Range("P188").Value = 73918.6
$ws.Range("Q188").Value = 110249.7
$ws.Range("R188").Value = 1301.5
$ws.Range("S188").Value = 50828.1
$ws.Range("T188").Value = 17869.1

# Row 189: new row for 22-09-2021
$ws.Range("A189").Value = "22-09-2021"
$ws.Range("B189").Value = 34258.3
$ws.Range("C189").Value = 14896.9
$ws.Range("D189").Value = 7083.4
$ws.Range("E189").Value = 29639.4
$ws.Range("F189").Value = 6637
$ws.Range("G189").Value = 15506.7
$ws.Range("I189").Value = 4821.8
$ws.Range("J189").Value = 1529
$ws.Range("K189").Value = 70950.8
$ws.Range("L189").Value = 4031
$ws.Range("M189").Value = 1619.6
$ws.Range("N189").Value = 16925.8
$ws.Range("O189").Value = 1407.5
$ws.Range("P189").Value = 74750.2
$ws.Range("Q189").Value = 112282.3
$ws.Range("R189").Value = 1310.3
$ws.Range("S189").Value = 51338.3
$ws.Range("T189").Value = 18099.1

# Row 190: new row for 23-09-2021
$ws.Range("A190").Value = "23-09-2021"
$ws.Range("B190").Value = 34764.8
$ws.Range("C190").Value = 15052.2
$ws.Range("D190").Value = 7078.4
$ws.Range("F190").Value = 6702
$ws.Range("G190").Value = 15644
$ws.Range("H190").Value = 3127.6
$ws.Range("I190").Value = 4853.2
$ws.Range("J190").Value = 1539.3
$ws.Range("K190").Value = 70725.5
$ws.Range("L190").Value = 4052.1
$ws.Range("M190").Value = 1631.2
$ws.Range("N190").Value = 17078.2
$ws.Range("O190").Value = 1401.5
$ws.Range("P190").Value = 74743.4
$ws.Range("Q190").Value = 114064.4
$ws.Range("R190").Value = 1312.2
$ws.Range("S190").Value = 51464.3
$ws.Range("T190").Value = 18208.3

# Row 191: new row for 24-09-2021
$ws.Range("A191").Value = "24-09-2021"
$ws.Range("B191").Value = 34798
$ws.Range("C191").Value = 15047.7
$ws.Range("D191").Value = 7051.5
$ws.Range("E191").Value = 30248.8
$ws.Range("F191").Value = 6638.5
$ws.Range("G191").Value = 15531.8
$ws.Range("H191").Value = 3125.2
$ws.Range("I191").Value = 4849.4
$ws.Range("J191").Value = 1532.1
$ws.Range("K191").Value = 70162.6
$ws.Range("L191").Value = 4038.2
$ws.Range("N191").Value = 17260.2
$ws.Range("O191").Value = 1384.7
$ws.Range("P191").Value = 74180.2
$ws.Range("Q191").Value = 113282.7
$ws.Range("R191").Value = 1309.4
$ws.Range("S191").Value = 51105.7
$ws.Range("T191").Value = 18308.9

# Row 192: new row for 27-09-2021
$ws.Range("A192").Value = "27-09-2021"
$ws.Range("D192").Value = 7065.6
$ws.Range("E192").Value = 30240.1
$ws.Range("F192").Value = 6668
$ws.Range("G192").Value = 15629.2
$ws.Range("H192").Value = 3133.6
$ws.Range("I192").Value = 4877.4
$ws.Range("J192").Value = 1533.1
$ws.Range("K192").Value = 70367.1
$ws.Range("M192").Value = 1620
$ws.Range("N192").Value = 17313.8
$ws.Range("O192").Value = 1395
